# Apply updated coin-ranking snapshot values (price/volume/hour refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '310.25'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.31%'
$ws.Range("E2").ClearFormats()
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '16'
$ws.Range("G2").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.37%'
$ws.Range("E3").ClearFormats()
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '16'
$ws.Range("G3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.198'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.39%'
$ws.Range("E4").ClearFormats()
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '16'
$ws.Range("G4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07696'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.61%'
$ws.Range("E5").ClearFormats()
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '16'
$ws.Range("G5").ClearFormats()
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.693'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '4.43%'
$ws.Range("E6").ClearFormats()
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '16'
$ws.Range("G6").ClearFormats()
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9457'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.95%'
$ws.Range("E7").ClearFormats()
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '16'
$ws.Range("G7").ClearFormats()
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.425'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.94%'
$ws.Range("E8").ClearFormats()
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '16'
$ws.Range("G8").ClearFormats()
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1266'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '8.73%'
$ws.Range("E9").ClearFormats()
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '16'
$ws.Range("G9").ClearFormats()
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1835'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.62%'
$ws.Range("E10").ClearFormats()
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '16'
$ws.Range("G10").ClearFormats()
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09114'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.93%'
$ws.Range("E11").ClearFormats()
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '16'
$ws.Range("G11").ClearFormats()
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04238'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.26%'
$ws.Range("E12").ClearFormats()
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '16'
$ws.Range("G12").ClearFormats()
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1053'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.93%'
$ws.Range("E13").ClearFormats()
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '16'
$ws.Range("G13").ClearFormats()
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001283'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '2.53%'
$ws.Range("E14").ClearFormats()
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '16'
$ws.Range("G14").ClearFormats()
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005892'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.55%'
$ws.Range("E15").ClearFormats()
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '16'
$ws.Range("G15").ClearFormats()
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.355'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.07%'
$ws.Range("E16").ClearFormats()
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '16'
$ws.Range("G16").ClearFormats()
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.295'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.38%'
$ws.Range("E17").ClearFormats()
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '16'
$ws.Range("G17").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '3.41%'
$ws.Range("E18").ClearFormats()
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '16'
$ws.Range("G18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.475'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '8.29%'
$ws.Range("E19").ClearFormats()
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '16'
$ws.Range("G19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1352'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.24%'
$ws.Range("E20").ClearFormats()
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '16'
$ws.Range("G20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2722'
$ws.Range("D21").ClearFormats()
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '16'
$ws.Range("G21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04028'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.03%'
$ws.Range("E22").ClearFormats()
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '16'
$ws.Range("G22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001267'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.73%'
$ws.Range("E23").ClearFormats()
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '16'
$ws.Range("G23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004239'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '3.40%'
$ws.Range("E24").ClearFormats()
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '16'
$ws.Range("G24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001272'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.12%'
$ws.Range("E25").ClearFormats()
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '16'
$ws.Range("G25").ClearFormats()
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '16'
$ws.Range("G26").ClearFormats()
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '16'
$ws.Range("G27").ClearFormats()
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '16'
$ws.Range("G28").ClearFormats()
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '16'
$ws.Range("G29").ClearFormats()
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '16'
$ws.Range("G30").ClearFormats()
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '16'
$ws.Range("G31").ClearFormats()
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '16'
$ws.Range("G32").ClearFormats()
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '16'
$ws.Range("G33").ClearFormats()
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '16'
$ws.Range("G34").ClearFormats()
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '16'
$ws.Range("G35").ClearFormats()
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '16'
$ws.Range("G36").ClearFormats()
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '16'
$ws.Range("G37").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '3.78%'
$ws.Range("E38").ClearFormats()
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '16'
$ws.Range("G38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05321'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.23%'
$ws.Range("E39").ClearFormats()
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '16'
$ws.Range("G39").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.65%'
$ws.Range("E40").ClearFormats()
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '16'
$ws.Range("G40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1315'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.96%'
$ws.Range("E41").ClearFormats()
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '16'
$ws.Range("G41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006920'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.86%'
$ws.Range("E42").ClearFormats()
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '16'
$ws.Range("G42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001943'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.39%'
$ws.Range("E43").ClearFormats()
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '16'
$ws.Range("G43").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-9.79%'
$ws.Range("E44").ClearFormats()
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '16'
$ws.Range("G44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3090'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.33%'
$ws.Range("E45").ClearFormats()
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '16'
$ws.Range("G45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006787'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.51%'
$ws.Range("E46").ClearFormats()
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '16'
$ws.Range("G46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.05%'
$ws.Range("E47").ClearFormats()
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '16'
$ws.Range("G47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2180'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '167.32%'
$ws.Range("E48").ClearFormats()
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '16'
$ws.Range("G48").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '3.36%'
$ws.Range("E49").ClearFormats()
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '16'
$ws.Range("G49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.05%'
$ws.Range("E50").ClearFormats()
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '16'
$ws.Range("G50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002003'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.05%'
$ws.Range("E51").ClearFormats()
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '16'
$ws.Range("G51").ClearFormats()
